# Roll the "legacy GSC export" date window forward:
#  - drop the oldest day (2025-10-20, currently row 2)
#  - append the 4 newest days (2026-01-16 .. 2026-01-19)
# Everything else (shared-string renumbering, the Table sheet's header
# refs, the dimension, and the HTTPS-URLs trend in column C) falls out
# naturally once the row shift happens.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Chart")

# 1) Drop the oldest dated row. This shifts every row below it up by
#    one, which is exactly why C4:C7 end up holding what used to be
#    C5:C8, etc.
$ws.Rows.Item(2).Delete()

# 2) Figure out where the data currently ends (row 88: header + 87 days)
#    and append the 4 new trailing days after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newDates = @("2026-01-16", "2026-01-17", "2026-01-18", "2026-01-19")

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = $lastRow + 1 + $i
    # Leading apostrophe forces text entry so "2026-01-16" isn't
    # reinterpreted as a date serial - matches the existing column A
    # cells, which store the dates as plain text.
    $ws.Cells.Item($r, 1).Value = "'" + $newDates[$i]
    $ws.Cells.Item($r, 2).Value = 0.0
    $ws.Cells.Item($r, 3).Value = 0.0
}

# 3) The apostrophe entry leaves a "quote prefix" text style on the new
#    A cells; re-stamp them with the plain (General) formatting used by
#    every other date cell in the column so the new rows match.
$ws.Range("A" + $lastRow).Copy()
$ws.Range("A" + ($lastRow + 1) + ":A" + ($lastRow + $newDates.Length)).PasteSpecial(-4122)
$excel.CutCopyMode = 0
